$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function SetText($sheet, $addr, $text) {
    $rng = $sheet.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

$ws.Range("D2").Value = '23.839.62'
$ws.Range("E2").Value = '  -3.18%  '
$ws.Range("D3").Value = '1.616.80'
$ws.Range("E3").Value = '  -3.60%  '
SetText $ws "D4" '0.9999'
$ws.Range("E4").Value = '  -0.07%  '
SetText $ws "D5" '307.28'
$ws.Range("E5").Value = '  -1.95%  '
SetText $ws "D6" '1.000'
$ws.Range("E6").Value = '  +0.04%  '
SetText $ws "D7" '0.3906'
$ws.Range("E7").Value = '  -0.60%  '
SetText $ws "D8" '0.3823'
$ws.Range("E8").Value = '  -3.42%  '
SetText $ws "D9" '0.9998'
$ws.Range("E9").Value = '  -0.10%  '
SetText $ws "D10" '1.363'
$ws.Range("E10").Value = '  -3.15%  '
SetText $ws "D11" '49.08'
$ws.Range("E11").Value = '  -3.65%  '
SetText $ws "D12" '0.08418'
$ws.Range("E12").Value = '  -2.91%  '
SetText $ws "D13" '23.91'
$ws.Range("E13").Value = '  -5.82%  '
SetText $ws "D14" '7.023'
$ws.Range("E14").Value = '  -4.40%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
SetText $ws "D15" '7.518'
$ws.Range("E15").Value = '  -2.72%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
SetText $ws "D16" '0.00001275'
$ws.Range("E16").Value = '  -3.65%  '
$ws.Range("D17").Value = '1.612.33'
$ws.Range("E17").Value = '  -2.95%  '
SetText $ws "D18" '93.11'
$ws.Range("E18").Value = '  -0.87%  '
SetText $ws "D19" '0.06908'
$ws.Range("E19").Value = '  -1.69%  '
SetText $ws "D20" '20.04'
$ws.Range("E20").Value = '  -7.07%  '
SetText $ws "D21" '6.812'
$ws.Range("E21").Value = '  -3.87%  '
SetText $ws "D22" '0.9999'
$ws.Range("E22").Value = '  +0.05%  '
SetText $ws "D23" '13.41'
$ws.Range("D24").Value = '23.848.96'
$ws.Range("E24").Value = '  -3.13%  '
SetText $ws "D25" '2.422'
$ws.Range("E25").Value = '  +2.35%  '
SetText $ws "D26" '2.838'
$ws.Range("E26").Value = '  +2.17%  '
SetText $ws "D27" '22.14'
$ws.Range("E27").Value = '  -4.42%  '
SetText $ws "D28" '157.12'
$ws.Range("E28").Value = '  -2.30%  '
SetText $ws "D29" '139.27'
$ws.Range("E29").Value = '  -5.19%  '
SetText $ws "D30" '5.249'
$ws.Range("E30").Value = '  -10.59%  '
SetText $ws "D31" '7.882'
$ws.Range("E31").Value = '  -5.46%  '
SetText $ws "D32" '2.482'
$ws.Range("E32").Value = '  -1.48%  '
$ws.Range("D33").Value = '1.790.45'
$ws.Range("E33").Value = '  -3.70%  '
SetText $ws "D34" '0.08077'
$ws.Range("E34").Value = '  -2.94%  '
SetText $ws "D35" '0.9754'
$ws.Range("E35").Value = '  -1.54%  '
SetText $ws "D36" '0.02872'
$ws.Range("E36").Value = '  -7.85%  '
$ws.Range("E37").Value = '  -5.51%  '
SetText $ws "D38" '0.2658'
$ws.Range("E38").Value = '  -5.61%  '
SetText $ws "D39" '0.09180'
$ws.Range("E39").Value = '  -3.66%  '
$ws.Range("E40").Value = '  -0.28%  '
SetText $ws "D41" '13.48'
$ws.Range("E41").Value = '  -1.09%  '
SetText $ws "D42" '1.431'
$ws.Range("E42").Value = '  -5.53%  '
SetText $ws "D43" '0.7454'
$ws.Range("E43").Value = '  -6.00%  '
SetText $ws "D44" '15.87'
$ws.Range("E44").Value = '  -4.43%  '
SetText $ws "D45" '0.6847'
$ws.Range("E45").Value = '  -4.17%  '
SetText $ws "D46" '2.453'
$ws.Range("E46").Value = '  -4.34%  '
$ws.Range("E47").Value = '  -2.68%  '
SetText $ws "D48" '1.000'
$ws.Range("E48").Value = '  +0.21%  '
SetText $ws "D49" '0.08258'
$ws.Range("E49").Value = '  -4.64%  '
SetText $ws "D50" '133.18'
$ws.Range("E50").Value = '  -3.25%  '
SetText $ws "D51" '1.209'
$ws.Range("E51").Value = '  -9.38%  '
